# "Finished Week 13 logging" - update Rushing & Receiving stats, add K.Allen
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Range("D2").Value = 5    # J.Herbert 2DATT
$rushing.Range("F2").Value = 8    # J.Herbert RZATT

$rushing.Range("C3").Value = 88   # A.Ekeler 1DATT
$rushing.Range("D3").Value = 50   # A.Ekeler 2DATT
$rushing.Range("F3").Value = 32   # A.Ekeler RZATT

$rushing.Range("C4").Value = 9    # J.Jackson 1DATT
$rushing.Range("D4").Value = 10   # J.Jackson 2DATT
$rushing.Range("E4").Value = 2    # J.Jackson 3DATT

$rushing.Range("C9").Value = 4    # J.Guyton 1DATT

# ---------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Insert a new row for K.Allen between L.Rountree (row 5) and
# M.Williams (previously row 6), then restore the bordered/bold
# style used by the rest of column A.
$receiving.Rows("6").Insert()
$receiving.Range("A5").Copy()
$receiving.Range("A6").PasteSpecial(-4122)

$receiving.Range("A6").Value = 4
$receiving.Range("B6").Value = "K.Allen"
$receiving.Range("C6").Value = 7
$receiving.Range("D6").Value = 5
$receiving.Range("E6").Value = 1
$receiving.Range("F6").Value = 0
$receiving.Range("G6").Value = 2
$receiving.Range("H6").Value = 2

# Existing player stat updates (rows reference the NEW row numbers,
# i.e. after the K.Allen row has shifted everything below it down by one)
$receiving.Range("C2").Value = 67   # A.Ekeler Short Target
$receiving.Range("D2").Value = 54   # A.Ekeler Short Comp

$receiving.Range("C3").Value = 10   # J.Jackson Short Target
$receiving.Range("D3").Value = 9    # J.Jackson Short Comp

$receiving.Range("C7").Value = 68   # M.Williams Short Target
$receiving.Range("D7").Value = 42   # M.Williams Short Comp
$receiving.Range("E7").Value = 25   # M.Williams Deep Target
$receiving.Range("F7").Value = 13   # M.Williams Deep Comp
$receiving.Range("G7").Value = 16   # M.Williams RZ Target

$receiving.Range("E8").Value = 2    # J.Palmer Deep Target

$receiving.Range("C9").Value = 24   # J.Guyton Short Target
$receiving.Range("D9").Value = 14   # J.Guyton Short Comp
$receiving.Range("E9").Value = 8    # J.Guyton Deep Target
$receiving.Range("F9").Value = 4    # J.Guyton Deep Comp

$receiving.Range("C11").Value = 47  # J.Cook Short Target
$receiving.Range("D11").Value = 30  # J.Cook Short Comp
$receiving.Range("E11").Value = 12  # J.Cook Deep Target
$receiving.Range("F11").Value = 6   # J.Cook Deep Comp

$receiving.Range("C12").Value = 22  # D.Parham Short Target
$receiving.Range("D12").Value = 17  # D.Parham Short Comp
$receiving.Range("E12").Value = 2   # D.Parham Deep Target
$receiving.Range("G12").Value = 5   # D.Parham RZ Target

# ---------------------------------------------------------------
# Finish on the Rushing tab (matches the saved workbook view)
# ---------------------------------------------------------------
$rushing.Activate()
$rushing.Range("A1").Select()
